{"js": "// Office.js (Word JS API) equivalent of the OOXML diff:\n// - The date heading paragraph text changes (\"2023-11-13 Monday\" -> \"2023-11-14 Tuesday\").\n// - Every one of the 100 arithmetic-problem cells in the single table gets a new\n//   \"a OP b=\" expression, in top-to-bottom / left-to-right (document) order.\n//\n// The document body, read top-to-bottom, is exactly: the heading paragraph followed\n// by the 100 table-cell paragraphs (one run each) \u2014 101 paragraphs total, which is\n// also the order `body.paragraphs` yields them in. We simply overwrite each\n// paragraph's text in place (insertText with \"Replace\") so the existing run\n// formatting (fonts/size) is preserved, and only the text content changes, matching\n// the diff exactly.\nconst newTexts = [\n  \"2023-11-14 Tuesday\", \"93-23=\", \"13+25=\", \"17+73=\", \"21+11=\", \"50-50=\", \"62-16=\", \"20+14=\",\n  \"55+15=\", \"11+55=\", \"26+50=\", \"24-14=\", \"46-43=\", \"58-4=\", \"68+3=\", \"21-2=\", \"97-64=\", \"55+21=\",\n  \"71+16=\", \"15+56=\", \"3+83=\", \"86-56=\", \"32+12=\", \"23+20=\", \"62-18=\", \"47+44=\", \"37+60=\", \"17+2=\",\n  \"64-49=\", \"67-18=\", \"90-71=\", \"30+31=\", \"57+30=\", \"2+59=\", \"57-23=\", \"77-75=\", \"87-32=\",\n  \"81-53=\", \"64+29=\", \"72+2=\", \"0+84=\", \"7+55=\", \"27+46=\", \"8+17=\", \"53-6=\", \"76-32=\", \"82+4=\",\n  \"42-31=\", \"83-41=\", \"48-8=\", \"22+63=\", \"58+20=\", \"4+87=\", \"49-40=\", \"55+25=\", \"27-8=\", \"23+69=\",\n  \"1+71=\", \"78-62=\", \"60-46=\", \"42-9=\", \"76+6=\", \"49+8=\", \"39+8=\", \"23-2=\", \"99-87=\", \"93-48=\",\n  \"93-12=\", \"38+4=\", \"0+95=\", \"66-59=\", \"84-14=\", \"86-7=\", \"95-43=\", \"79-37=\", \"91-79=\", \"61-51=\",\n  \"52+8=\", \"74-7=\", \"79-38=\", \"1+23=\", \"61-43=\", \"19+68=\", \"84-62=\", \"84+10=\", \"97-70=\", \"81-18=\",\n  \"84-82=\", \"31-28=\", \"84-51=\", \"29-14=\", \"95-64=\", \"82-63=\", \"18-0=\", \"45+32=\", \"24+16=\", \"71+9=\",\n  \"71-7=\", \"97-42=\", \"0+40=\", \"72-4=\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newTexts.length) {\n  throw new Error(\n    \"Paragraph count mismatch: expected \" + newTexts.length + \" got \" + paragraphs.items.length\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].insertText(newTexts[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM equivalent of the OOXML diff:\n# - The date heading paragraph text changes (\"2023-11-13 Monday\" -> \"2023-11-14 Tuesday\").\n# - Every one of the 100 arithmetic-problem cells in the single table gets a new\n#   \"a OP b=\" expression, in top-to-bottom / left-to-right (row-major / document) order.\n#\n# Each cell/paragraph Range.Text ends with a hidden mark (paragraph mark \"\\r\" for the\n# heading, cell-end mark for table cells) that must be excluded before assigning new\n# text, otherwise it would be clobbered. MoveEnd(wdCharacter, -1) backs off exactly\n# that one trailing mark as a single navigable unit, leaving only the visible text\n# selected; assigning .Text then swaps the content while keeping the run formatting\n# (fonts/size) intact, matching the diff exactly.\n\n$d = $word.ActiveDocument\n\n# Update the date heading (first paragraph, outside the table).\n$titlePara = $d.Paragraphs.Item(1)\n$titleRange = $titlePara.Range\n[void]$titleRange.MoveEnd(1, -1)  # trim the paragraph mark\n$titleRange.Text = \"2023-11-14 Tuesday\"\n\n# Update every cell of the first (only) table, left-to-right, top-to-bottom,\n# matching document order.\n$newValues = @(\n    \"93-23=\",\n    \"13+25=\",\n    \"17+73=\",\n    \"21+11=\",\n    \"50-50=\",\n    \"62-16=\",\n    \"20+14=\",\n    \"55+15=\",\n    \"11+55=\",\n    \"26+50=\",\n    \"24-14=\",\n    \"46-43=\",\n    \"58-4=\",\n    \"68+3=\",\n    \"21-2=\",\n    \"97-64=\",\n    \"55+21=\",\n    \"71+16=\",\n    \"15+56=\",\n    \"3+83=\",\n    \"86-56=\",\n    \"32+12=\",\n    \"23+20=\",\n    \"62-18=\",\n    \"47+44=\",\n    \"37+60=\",\n    \"17+2=\",\n    \"64-49=\",\n    \"67-18=\",\n    \"90-71=\",\n    \"30+31=\",\n    \"57+30=\",\n    \"2+59=\",\n    \"57-23=\",\n    \"77-75=\",\n    \"87-32=\",\n    \"81-53=\",\n    \"64+29=\",\n    \"72+2=\",\n    \"0+84=\",\n    \"7+55=\",\n    \"27+46=\",\n    \"8+17=\",\n    \"53-6=\",\n    \"76-32=\",\n    \"82+4=\",\n    \"42-31=\",\n    \"83-41=\",\n    \"48-8=\",\n    \"22+63=\",\n    \"58+20=\",\n    \"4+87=\",\n    \"49-40=\",\n    \"55+25=\",\n    \"27-8=\",\n    \"23+69=\",\n    \"1+71=\",\n    \"78-62=\",\n    \"60-46=\",\n    \"42-9=\",\n    \"76+6=\",\n    \"49+8=\",\n    \"39+8=\",\n    \"23-2=\",\n    \"99-87=\",\n    \"93-48=\",\n    \"93-12=\",\n    \"38+4=\",\n    \"0+95=\",\n    \"66-59=\",\n    \"84-14=\",\n    \"86-7=\",\n    \"95-43=\",\n    \"79-37=\",\n    \"91-79=\",\n    \"61-51=\",\n    \"52+8=\",\n    \"74-7=\",\n    \"79-38=\",\n    \"1+23=\",\n    \"61-43=\",\n    \"19+68=\",\n    \"84-62=\",\n    \"84+10=\",\n    \"97-70=\",\n    \"81-18=\",\n    \"84-82=\",\n    \"31-28=\",\n    \"84-51=\",\n    \"29-14=\",\n    \"95-64=\",\n    \"82-63=\",\n    \"18-0=\",\n    \"45+32=\",\n    \"24+16=\",\n    \"71+9=\",\n    \"71-7=\",\n    \"97-42=\",\n    \"0+40=\",\n    \"72-4=\"\n)\n\n$tbl = $d.Tables.Item(1)\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\n\nif (($rows * $cols) -ne $newValues.Count) {\n    throw \"Cell count mismatch: expected $($newValues.Count) got $($rows * $cols)\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cellRange = $cell.Range\n        [void]$cellRange.MoveEnd(1, -1)  # trim the cell-end mark\n        $cellRange.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
